$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    12033, 12033, 11706, 11706, 11706, 10837, 10837, 10837, 10837, 10837, 9718, 9718, 9553, 9553, 9553, 9553, 9265, 9265, 9126, 9126,
    9126, 9126, 9126, 9126, 9126, 9126, 9126, 9126, 8842, 8842, 8842, 8842, 8842, 8842, 8842, 8842, 8842, 8842, 8842, 8842,
    8842, 8842, 8842, 8762, 8762, 8762, 8762, 8700, 8700, 8700, 8700, 8700, 8700, 8700, 8700, 8660, 8502, 8362, 8362, 8362,
    8362, 8362, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025,
    8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 8025, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569,
    7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569, 7569
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}
